$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the header formatting used by the rest of row 1 by copying H1's
# format (bold, thin box border, centered horizontally, top-aligned
# vertically) over to the new I1:J1 header cells.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 5
